# Auto-generated edit script: updates market-price columns (H-N) across
# all 8 Leve-profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to match
# the latest scheduled-runner price snapshot.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 106.583336
$ws.Range("H17").Value = 37594.71
$ws.Range("J17").Value = 38286.277
$ws.Range("L17").Value = 114858.831
$ws.Range("N17").Value = -115194.831
$ws.Range("H18").Value = 1439
$ws.Range("J18").Value = 2997.5
$ws.Range("L18").Value = 2997.5
$ws.Range("N18").Value = -3565.5
$ws.Range("H19").Value = 3226.2
$ws.Range("I19").Value = 573.5
$ws.Range("K19").Value = 573.5
$ws.Range("M19").Value = -398.5
$ws.Range("H51").Value = 9379.799999999999
$ws.Range("I51").Value = 8949.5
$ws.Range("K51").Value = 8949.5
$ws.Range("M51").Value = -8465.5
$ws.Range("H69").Value = 11751.6
$ws.Range("I69").Value = 8586.333000000001
$ws.Range("J69").Value = 16499.5
$ws.Range("K69").Value = 25758.999
$ws.Range("L69").Value = 49498.5
$ws.Range("M69").Value = -24884.999
$ws.Range("N69").Value = -51246.5
$ws.Range("H72").Value = 11751.6
$ws.Range("I72").Value = 8586.333000000001
$ws.Range("J72").Value = 16499.5
$ws.Range("K72").Value = 77276.997
$ws.Range("L72").Value = 148495.5
$ws.Range("M72").Value = -72908.997
$ws.Range("N72").Value = -157231.5
$ws.Range("H101").Value = 678
$ws.Range("J101").Value = 1000
$ws.Range("L101").Value = 3000
$ws.Range("N101").Value = -6244
$ws.Range("H127").Value = 20303.889
$ws.Range("I127").Value = 21741.875
$ws.Range("K127").Value = 65225.625
$ws.Range("M127").Value = -60265.625
$ws.Range("H130").Value = 80000
$ws.Range("J130").Value = 80000
$ws.Range("L130").Value = 80000
$ws.Range("N130").Value = -90040
$ws.Range("H132").Value = 2301.853
$ws.Range("I132").Value = 1402.25
$ws.Range("K132").Value = 4206.75
$ws.Range("M132").Value = -1676.75
$ws.Range("H133").Value = 62859.547
$ws.Range("J133").Value = 62859.547
$ws.Range("L133").Value = 62859.547
$ws.Range("N133").Value = -72979.54699999999
$ws.Range("H134").Value = 75000
$ws.Range("J134").Value = 75000
$ws.Range("L134").Value = 75000
$ws.Range("N134").Value = -85140
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()
$ws.Range("H137").Value = 3736
$ws.Range("I137").Value = 1903.6471
$ws.Range("J137").Value = 5151.909
$ws.Range("K137").Value = 5710.9413
$ws.Range("L137").Value = 15455.727
$ws.Range("M137").Value = -3160.9413
$ws.Range("N137").Value = -20555.727
$ws.Range("H138").Value = 4450
$ws.Range("I138").Value = 1678.4445
$ws.Range("K138").Value = 5035.333500000001
$ws.Range("M138").Value = 104.6664999999994

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 200737.8
$ws.Range("I2").Value = 250672.25
$ws.Range("K2").Value = 250672.25
$ws.Range("M2").Value = -250559.25
$ws.Range("H32").Value = 3666.3115
$ws.Range("I32").Value = 2998.2
$ws.Range("K32").Value = 2998.2
$ws.Range("M32").Value = -2711.2
$ws.Range("H45").Value = 7997.3076
$ws.Range("I45").Value = 3796.1667
$ws.Range("J45").Value = 11598.286
$ws.Range("K45").Value = 3796.1667
$ws.Range("L45").Value = 11598.286
$ws.Range("M45").Value = -3419.1667
$ws.Range("N45").Value = -12352.286
$ws.Range("H46").Value = 5713
$ws.Range("I46").Value = 4999
$ws.Range("J46").Value = 6070
$ws.Range("K46").Value = 4999
$ws.Range("L46").Value = 6070
$ws.Range("M46").Value = -4680
$ws.Range("N46").Value = -6708
$ws.Range("H74").Value = 2068.258
$ws.Range("I74").Value = 1452.5555
$ws.Range("K74").Value = 1452.5555
$ws.Range("M74").Value = -578.5554999999999
$ws.Range("H77").Value = 2068.258
$ws.Range("I77").Value = 1452.5555
$ws.Range("K77").Value = 7262.7775
$ws.Range("M77").Value = -2894.7775
$ws.Range("H116").Value = 200737.8
$ws.Range("I116").Value = 250672.25
$ws.Range("K116").Value = 250672.25
$ws.Range("M116").Value = -248378.25
$ws.Range("H132").Value = 3254.3333
$ws.Range("I132").Value = 899.3
$ws.Range("K132").Value = 2697.9
$ws.Range("M132").Value = -167.8999999999996

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 200737.8
$ws.Range("I3").Value = 250672.25
$ws.Range("K3").Value = 250672.25
$ws.Range("M3").Value = -250558.25
$ws.Range("H7").Value = 716.6667
$ws.Range("I7").Value = 75
$ws.Range("K7").Value = 75
$ws.Range("M7").Value = 38
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H86").Value = 774264.3
$ws.Range("I86").Value = 811000.4
$ws.Range("K86").Value = 811000.4
$ws.Range("M86").Value = -809877.4
$ws.Range("H89").Value = 774264.3
$ws.Range("I89").Value = 811000.4
$ws.Range("K89").Value = 4055002
$ws.Range("M89").Value = -4049386
$ws.Range("H134").Value = 4004.3062
$ws.Range("I134").Value = 2803.7646
$ws.Range("J134").Value = 6725.533
$ws.Range("K134").Value = 8411.293799999999
$ws.Range("L134").Value = 20176.599
$ws.Range("M134").Value = -5876.293799999999
$ws.Range("N134").Value = -25246.599

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 287
$ws.Range("I22").Value = 287
$ws.Range("K22").Value = 287
$ws.Range("M22").Value = 63
$ws.Range("H31").Value = 2586.158
$ws.Range("I31").Value = 1490.88
$ws.Range("J31").Value = 3441.8438
$ws.Range("K31").Value = 1490.88
$ws.Range("L31").Value = 3441.8438
$ws.Range("M31").Value = -1195.88
$ws.Range("N31").Value = -4031.8438
$ws.Range("H34").Value = 2586.158
$ws.Range("I34").Value = 1490.88
$ws.Range("J34").Value = 3441.8438
$ws.Range("K34").Value = 1490.88
$ws.Range("L34").Value = 3441.8438
$ws.Range("M34").Value = -1288.88
$ws.Range("N34").Value = -3845.8438
$ws.Range("H53").Value = 79999
$ws.Range("J53").Value = 79999
$ws.Range("L53").Value = 79999
$ws.Range("N53").Value = -81213

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 6666746.5
$ws.Range("I11").Value = 10000000
$ws.Range("J11").Value = 240
$ws.Range("K11").Value = 30000000
$ws.Range("L11").Value = 720
$ws.Range("M11").Value = -29999860
$ws.Range("N11").Value = -1000
$ws.Range("H127").Value = 1962.6666
$ws.Range("J127").Value = 1962.6666
$ws.Range("L127").Value = 5887.9998
$ws.Range("N127").Value = -15807.9998
$ws.Range("H139").Value = 1882.9166
$ws.Range("I139").Value = 1882.9166
$ws.Range("K139").Value = 5648.7498
$ws.Range("M139").Value = -508.7497999999996

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H110").Value = 50000
$ws.Range("J110").Value = 50000
$ws.Range("L110").Value = 50000
$ws.Range("N110").Value = -58180
$ws.Range("H135").Value = 112499.5
$ws.Range("J135").Value = 112499.5
$ws.Range("L135").Value = 112499.5
$ws.Range("N135").Value = -122639.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7254.778
$ws.Range("I7").Value = 7553.273
$ws.Range("J7").Value = 6785.7144
$ws.Range("K7").Value = 7553.273
$ws.Range("L7").Value = 6785.7144
$ws.Range("M7").Value = -7441.273
$ws.Range("N7").Value = -7009.7144
$ws.Range("H16").Value = 380
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 380
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 380
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -720
$ws.Range("H55").Value = 1094
$ws.Range("I55").Value = 410.2857
$ws.Range("K55").Value = 410.2857
$ws.Range("M55").Value = -237.2857
$ws.Range("H122").Value = 748721.5
$ws.Range("I122").Value = 725442.6
$ws.Range("J122").Value = 773791.0600000001
$ws.Range("K122").Value = 2176327.8
$ws.Range("L122").Value = 2321373.18
$ws.Range("M122").Value = -2173877.8
$ws.Range("N122").Value = -2326273.18
$ws.Range("H126").Value = 7254.778
$ws.Range("I126").Value = 7553.273
$ws.Range("J126").Value = 6785.7144
$ws.Range("K126").Value = 22659.819
$ws.Range("L126").Value = 20357.1432
$ws.Range("M126").Value = -20189.819
$ws.Range("N126").Value = -25297.1432
$ws.Range("H132").Value = 4090.125
$ws.Range("I132").Value = 3214.95
$ws.Range("K132").Value = 9644.849999999999
$ws.Range("M132").Value = -7114.849999999999
$ws.Range("H136").Value = 3698.5715
$ws.Range("I136").Value = 2267.5334
$ws.Range("J136").Value = 4771.85
$ws.Range("K136").Value = 6802.600199999999
$ws.Range("L136").Value = 14315.55
$ws.Range("M136").Value = -4252.600199999999
$ws.Range("N136").Value = -19415.55

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 14068.143
$ws.Range("J45").Value = 13595.6
$ws.Range("L45").Value = 13595.6
$ws.Range("N45").Value = -14577.6
$ws.Range("H126").Value = 3162.6155
$ws.Range("I126").Value = 3219.8572
$ws.Range("K126").Value = 9659.571599999999
$ws.Range("M126").Value = -7189.571599999999
$ws.Range("H132").Value = 2830.9656
$ws.Range("I132").Value = 1878.2
$ws.Range("K132").Value = 5634.6
$ws.Range("M132").Value = -3104.6
